$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.640.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.267.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.58"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.68%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.608.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.99"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.266.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.544.61"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0993"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.02"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.71"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.89"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0853"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.19"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.86"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +12.11%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.100"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.30%  "
